# Update MRVL balance sheet figures (BalanceSheet/MRVL_bal.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MRVL")

# Row 4 - Inventory
$ws.Range("B4").Value = 268000000.0
$ws.Range("C4").Value = 268000000.0
$ws.Range("D4").Value = 263000000.0
$ws.Range("E4").Value = 270000000.0
$ws.Range("F4").Value = 323000000.0

# Row 14 - Accounts Payable
$ws.Range("B14").Value = 252000000.0
$ws.Range("C14").Value = 224000000.0
$ws.Range("D14").Value = 239000000.0
$ws.Range("E14").Value = 186000000.0
$ws.Range("F14").Value = 214000000.0

# Row 21 - Long Term Tax Liability (Deferred)
$ws.Range("B21").Value = -650000000.0
$ws.Range("C21").Value = -609000000.0
$ws.Range("D21").Value = -609000000.0
$ws.Range("E21").Value = -606000000.0
$ws.Range("F21").Value = -609000000.0
